$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to remain plain text so values like "30.391.30" or
# trailing-zero decimals ("0.9990") are not auto-coerced into numbers.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '30.391.30'
$ws.Range("E2").Value = '  -0.03%  '

# Row 3
$ws.Range("D3").Value = '1.885.33'
$ws.Range("E3").Value = '  -1.07%  '

# Row 4
$ws.Range("D4").Value = '0.9997'
$ws.Range("E4").Value = '  +0.02%  '

# Row 5
$ws.Range("D5").Value = '238.28'
$ws.Range("E5").Value = '  -0.13%  '

# Row 6
$ws.Range("D6").Value = '0.9994'
$ws.Range("E6").Value = '  +0.04%  '

# Row 7
$ws.Range("D7").Value = '0.4685'
$ws.Range("E7").Value = '  -0.85%  '

# Row 8
$ws.Range("D8").Value = '0.2826'
$ws.Range("E8").Value = '  -0.25%  '

# Row 9
$ws.Range("D9").Value = '0.06577'
$ws.Range("E9").Value = '  -1.12%  '

# Row 10
$ws.Range("D10").Value = '19.77'
$ws.Range("E10").Value = '  +5.90%  '

# Row 11
$ws.Range("E11").Value = '  -1.70%  '

# Row 12
$ws.Range("D12").Value = '0.07759'
$ws.Range("E12").Value = '  +0.54%  '

# Row 13
$ws.Range("D13").Value = '1.881.61'
$ws.Range("E13").Value = '  -1.18%  '

# Row 14
$ws.Range("D14").Value = '5.138'
$ws.Range("E14").Value = '  -1.00%  '

# Row 15
$ws.Range("D15").Value = '0.6681'
$ws.Range("E15").Value = '  -0.14%  '

# Row 16
$ws.Range("D16").Value = '285.82'
$ws.Range("E16").Value = '  +12.64%  '

# Row 17
$ws.Range("D17").Value = '30.390.50'
$ws.Range("E17").Value = '  +0.04%  '

# Row 18
$ws.Range("D18").Value = '0.9990'
$ws.Range("E18").Value = '  -0.07%  '

# Row 19
$ws.Range("D19").Value = '12.65'
$ws.Range("E19").Value = '  +0.37%  '

# Row 20
$ws.Range("D20").Value = '2.129.42'
$ws.Range("E20").Value = '  -1.18%  '

# Row 21
$ws.Range("B21").Value = 'Uniswap'
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D21").Value = '5.375'
$ws.Range("E21").Value = '  +0.29%  '

# Row 22
$ws.Range("B22").Value = 'ShibaInu'
$ws.Range("C22").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D22").Value = '0.000007310'
$ws.Range("E22").Value = '  -1.68%  '

# Row 23
$ws.Range("D23").Value = '0.9995'
$ws.Range("E23").Value = '  +0.06%  '

# Row 24
$ws.Range("B24").Value = 'Chainlink'
$ws.Range("C24").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D24").Value = '6.190'
$ws.Range("E24").Value = '  -1.50%  '

# Row 25
$ws.Range("B25").Value = 'Monero'
$ws.Range("C25").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D25").Value = '167.35'
$ws.Range("E25").Value = '  +0.10%  '

# Row 26
$ws.Range("D26").Value = '9.280'
$ws.Range("E26").Value = '  -0.47%  '

# Row 27
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").Value = '19.11'
$ws.Range("E27").Value = '  +1.24%  '

# Row 28
$ws.Range("B28").Value = 'LidoDAOToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D28").Value = '1.985'
$ws.Range("E28").Value = '  -2.99%  '

# Row 29
$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").Value = '1.374'
$ws.Range("E29").Value = '  -0.12%  '

# Row 30
$ws.Range("B30").Value = 'Stellar'
$ws.Range("C30").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D30").Value = '0.09839'
$ws.Range("E30").Value = '  -2.46%  '

# Row 31
$ws.Range("B31").Value = 'Filecoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D31").Value = '4.475'
$ws.Range("E31").Value = '  -4.46%  '

# Row 32
$ws.Range("B32").Value = 'PancakeSwap'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D32").Value = '1.494'
$ws.Range("E32").Value = '  -0.92%  '

# Row 33
$ws.Range("B33").Value = 'InternetComputer(DFINITY)'
$ws.Range("C33").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D33").Value = '4.188'
$ws.Range("E33").Value = '  -1.78%  '

# Row 34
$ws.Range("B34").Value = 'Hedera'
$ws.Range("C34").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D34").Value = '0.04719'
$ws.Range("E34").Value = '  +0.24%  '

# Row 35
$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D35").Value = '0.7110'
$ws.Range("E35").Value = '  -2.36%  '

# Row 36
$ws.Range("B36").Value = 'ARBITRUM'
$ws.Range("C36").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D36").Value = '1.098'
$ws.Range("E36").Value = '  -0.86%  '

# Row 37
$ws.Range("B37").Value = 'HuobiToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D37").Value = '2.717'
$ws.Range("E37").Value = '  +0.66%  '

# Row 38
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").Value = '0.01875'
$ws.Range("E38").Value = '  -2.05%  '

# Row 39
$ws.Range("B39").Value = 'FraxShare'
$ws.Range("C39").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D39").Value = '6.745'
$ws.Range("E39").Value = '  +8.33%  '

# Row 40
$ws.Range("B40").Value = 'MXToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D40").Value = '2.537'
$ws.Range("E40").Value = '  -2.06%  '

# Row 41
$ws.Range("B41").Value = 'Aave'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D41").Value = '72.79'
$ws.Range("E41").Value = '  +0.57%  '

# Row 42
$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").Value = '0.8785'
$ws.Range("E42").Value = '  +2.90%  '

# Row 43
$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D43").Value = '1.972'
$ws.Range("E43").Value = '  +0.53%  '

# Row 44
$ws.Range("B44").Value = 'Quant'
$ws.Range("C44").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D44").Value = '104.57'
$ws.Range("E44").Value = '  -1.34%  '

# Row 45
$ws.Range("B45").Value = 'PaxDollar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D45").Value = '0.9989'
$ws.Range("E45").Value = '  +0.02%  '

# Row 46
$ws.Range("B46").Value = 'TheSandbox'
$ws.Range("C46").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D46").Value = '0.4203'
$ws.Range("E46").Value = '  -0.54%  '

# Row 47
$ws.Range("B47").Value = 'Maker'
$ws.Range("C47").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D47").Value = '992.35'
$ws.Range("E47").Value = '  -0.85%  '

# Row 48
$ws.Range("B48").Value = 'Aptos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D48").Value = '7.258'
$ws.Range("E48").Value = '  -1.87%  '

# Row 49
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").Value = '9.372'
$ws.Range("E49").Value = '  +7.29%  '

# Row 50
$ws.Range("B50").Value = 'Elrond'
$ws.Range("C50").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D50").Value = '34.24'
$ws.Range("E50").Value = '  -0.63%  '

# Row 51
$ws.Range("B51").Value = 'Algorand'
$ws.Range("C51").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D51").Value = '0.1164'
$ws.Range("E51").Value = '  -2.39%  '
